$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Add header for new column D
$ws.Range("D1").Value = "percent_of_control"

# Add formulas for D2:D6 (percent of control, relative to I1)
$ws.Range("D2").Formula = "=(B2/`$I`$1)*100"
$ws.Range("D3:D6").Formula = "=(B3/`$I`$1)*100"

# Update the selected cell to D6, matching the author's last-edited cell
$ws.Activate()
$ws.Range("D6").Select()
